$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154 -- this shifts the existing rows 154:163 down to
# 155:164 (matching the diff, which is a weekly-data insert at the top of
# this Berenjena / Macroferia Regional de Talca block) and grows the used
# range to A1:R164.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with this week's record. The
# "constant" columns (A,B,C,E,F,G,H,I,R) carry the same values as every
# other row in this block.
$ws.Cells.Item(154, 1).Value = 5
$ws.Cells.Item(154, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(154, 3).Value = "Maule"
$ws.Range("D154").Value = 45021
$ws.Cells.Item(154, 5).Value = 7
$ws.Cells.Item(154, 6).Value = 100112001
$ws.Cells.Item(154, 7).Value = "Berenjena"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 150
$ws.Cells.Item(154, 11).Value = 8000
$ws.Cells.Item(154, 12).Value = 8000
$ws.Cells.Item(154, 13).Value = 8000
$ws.Cells.Item(154, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(154, 15).Value = "Región del Maule"
$ws.Cells.Item(154, 16).Value = 160
$ws.Cells.Item(154, 17).Value = 50
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# Apply the same date-time number format (style index 2 in before.xlsx) the
# rest of column D uses, so the new row matches its neighbours.
$ws.Range("D154").NumberFormat = $ws.Range("D155").NumberFormat
